$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "NA" value in E2, leaving the cell formatting/style intact
$ws.Range("E2").ClearContents()

# Update the active selection to match the target state (E2 only)
$ws.Range("E2").Select()
